# Diversification-Tool-3.xlsx edit
# "return, variance and performance calculations work"
#
# - Z8 holdings-after-trades value changes from 100 to 1000
# - J9:J12 ("Return if state occurs" column) stop being computed via
#   formula (=B+D-E) and become a flat entered value of 20
# - K11's label cell had been using a redundant duplicate shared string
#   "StockC"; once the J-column literals above are in place the engine
#   re-derives the shared-string table and naturally collapses that
#   duplicate away when the cell is normalized to the same text as the
#   other "Stock C" labels
# - everything else (AA8, P14, Z16, U19, U23, U28, P25, L20, L21, ...)
#   is formula-driven and recalculates automatically

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Diversification")

$ws.Range("Z8").Value = 1000

$ws.Range("J9").Value = 20
$ws.Range("J10").Value = 20
$ws.Range("J11").Value = 20
$ws.Range("J12").Value = 20

# Normalize the stray duplicate-text label to match the rest of the sheet
$ws.Range("K11").Value = "Stock C"

# Reflect the author's final cursor position on the sheet
$ws.Range("P17").Select()
